# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.502.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.615.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.71%  "
$ws.Range("E7").Value = "  +3.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.609.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -5.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +13.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.616"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000281"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.198.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "674.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.621.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.583.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.122"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.936"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.32%  "
$ws.Range("E26").Value = "  -3.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("E34").Value = "  -5.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.93"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "569.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.107"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.26%  "
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0450"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.95%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.528.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.12%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.345"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("E44").Value = "  -3.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "34.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0726"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.10%  "
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.56%  "
